$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costs")

# Rename "battery_link" (row 9) to "batterylink"
$ws.Range("A9").Value = "batterylink"

# Update selection to A10
$ws.Range("A10").Select()
